$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4581
$ws.Range("L3").Value = 4906
$ws.Range("L4").Value = 1216
$ws.Range("L6").Value = 4194
$ws.Range("L7").Value = 15180

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L7").Value = 501
$ws.Range("L14").Value = 76
$ws.Range("L19").Value = 418
$ws.Range("L20").Value = 387
$ws.Range("L23").Value = 164
$ws.Range("L29").Value = 842
$ws.Range("L31").Value = 147
$ws.Range("L33").Value = 688
$ws.Range("L36").Value = 197
$ws.Range("L37").Value = 559
$ws.Range("L44").Value = 110
$ws.Range("L45").Value = 27
$ws.Range("L48").Value = 196
$ws.Range("L51").Value = 185
$ws.Range("L53").Value = 175
$ws.Range("L54").Value = 317
$ws.Range("L59").Value = 28
$ws.Range("L63").Value = 47
$ws.Range("L64").Value = 104
$ws.Range("L65").Value = 294
$ws.Range("L66").Value = 36
$ws.Range("L67").Value = 521
$ws.Range("L72").Value = 60
$ws.Range("L76").Value = 237
$ws.Range("L79").Value = 401
$ws.Range("L83").Value = 330
$ws.Range("L85").Value = 778
$ws.Range("L89").Value = 217
$ws.Range("L90").Value = 151
$ws.Range("L95").Value = 205
$ws.Range("L98").Value = 85
$ws.Range("L99").Value = 260
$ws.Range("L101").Value = 15180

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 169
$ws.Range("L3").Value = 164
$ws.Range("L7").Value = 501

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 60
$ws.Range("L4").Value = 35
$ws.Range("L7").Value = 217

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 235
$ws.Range("L7").Value = 778

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 129
$ws.Range("L7").Value = 330

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L6").Value = 213
$ws.Range("L7").Value = 688

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 65
$ws.Range("L6").Value = 46
$ws.Range("L7").Value = 205

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 194
$ws.Range("L7").Value = 559

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 105
$ws.Range("L3").Value = 92
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 294

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 70
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 58
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L6").Value = 118
$ws.Range("L7").Value = 521

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L2").Value = 59
$ws.Range("L7").Value = 317

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 253
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 842

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 418

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L6").Value = 85
$ws.Range("L7").Value = 401

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 116
$ws.Range("L7").Value = 387

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 197

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 185

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 27
